$wb = $excel.ActiveWorkbook
$wsProperty = $wb.Worksheets.Item("Property")
$wsProperty.Range("C20").Value = "//div[(@id='gridProperty')]"
$wsProperty.Range("A21").Value = "tbLease"
$wsProperty.Range("B21").Value = "by_xpath"
$wsProperty.Range("C21").Value = "//div[(@id='id_376368_C_ctl01_ucPSPC_gvPropertyContracts')]"

$wsSearch = $wb.Worksheets.Item("GlobalSearch")
$wsSearch.Rows(4).Insert()
$wsSearch.Range("A4").Value = "searchOption"
$wsSearch.Range("B4").Value = "by_xpath"
$wsSearch.Range("C4").Value = "//*[contains(@id,'ddlSearchOptions')]"
$wsSearch.Range("D4").Font.Name = "Consolas"
$wsSearch.Range("D4").Font.Size = 9
$wsSearch.Range("D4").Font.Family = 3
$wsSearch.Range("D4").Font.Color = 2236962

# Restore view/selection state to match the saved workbook: GlobalSearch
# ends up with D4 selected, while Property stays the active/selected tab
# with C23 selected.
$wsSearch.Range("D4").Select()
$wsProperty.Activate()
$wsProperty.Range("C23").Select()
